# Rename the first sheet ("INTER_SWITCH_LINKS" -> "SWITCH_TO_SWITCH") and
# move the active/selected sheet + selection from HARDWARE_MANAGEMENT
# (previously active, cell D28 selected) to the renamed sheet (now active,
# with F41 selected).

$wb = $excel.ActiveWorkbook

# 1. Rename INTER_SWITCH_LINKS -> SWITCH_TO_SWITCH
$switchSheet = $wb.Worksheets.Item("INTER_SWITCH_LINKS")
$switchSheet.Name = "SWITCH_TO_SWITCH"

# 2. Make it the active sheet and select F41 on it (new tabSelected / selection)
$switchSheet.Activate()
$switchSheet.Range("F41").Select()
